$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Add the new row of data (row 15)
$ws.Range("A15").Value = "ip_banned"
$ws.Range("B15").Value = "Insert"
$ws.Range("C15").Value = "REALMS_INS_BANIP_BANAUTOIP"
$ws.Range("D15").Value = "Ban auto de l'ip pour erreur authentification"
$ws.Range("E15").Value = "INSERT INTO ip_banned VALUES (`$1, now(), now() + INTERVAL  '20 minute', 'AutoBan pour erreur authentification ', `$2, true)"

# Update the view: scroll so column C is the leftmost visible column, and select D16
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("D16").Select()
